$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.885.10"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.752.17"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.98"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.50"
$ws.Range("E6").Value = "  +4.17%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.68"
$ws.Range("E11").Value = "  -15.66%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "3.239.94"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.01"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "63.850.48"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "2.753.54"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.40"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.569"
$ws.Range("E22").Value = "  +5.81%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.40"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.68"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "0.0₃0940"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  +4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.35"
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.60"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.01"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("E40").Value = "  +8.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "332.49"
$ws.Range("E41").Value = "  -4.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.53"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.11"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0600"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.04"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "137.11"
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  +0.81%  "
